# Insert a floating "Straight Arrow Connector" shape (wrapped in an
# mc:AlternateContent run, with VML fallback) as the first run of the
# document's only paragraph, right before the existing w:bookmarkStart.
#
# Real Word exposes no Shapes.AddConnector-equivalent that lands inside a
# specific run position reliably across hosts, so we build the exact OOXML
# for the run (mirroring what Word itself emits when you draw a Straight
# Arrow Connector) and splice it in via Range.InsertXML on a Range collapsed
# to the start of the target paragraph - this merges the new run into that
# paragraph instead of creating a sibling paragraph.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$insertionRange = $p1.Range
$insertionRange.Collapse(1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>5463540</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>685800</wp:posOffset></wp:positionV><wp:extent cx="0" cy="327660"/><wp:effectExtent l="95250" t="19050" r="114300" b="91440"/><wp:wrapNone/><wp:docPr id="2" name="Straight Arrow Connector 2"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvCnPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="0" cy="327660"/></a:xfrm><a:prstGeom prst="straightConnector1"><a:avLst/></a:prstGeom><a:ln><a:tailEnd type="arrow"/></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="dk1"/></a:lnRef><a:fillRef idx="0"><a:schemeClr val="dk1"/></a:fillRef><a:effectRef idx="1"><a:schemeClr val="dk1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="tx1"/></a:fontRef></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype id="_x0000_t32" coordsize="21600,21600" o:spt="32" o:oned="t" path="m,l21600,21600e" filled="f"><v:path arrowok="t" fillok="f" o:connecttype="none"/><o:lock v:ext="edit" shapetype="t"/></v:shapetype><v:shape id="Straight Arrow Connector 2" o:spid="_x0000_s1026" type="#_x0000_t32" style="position:absolute;margin-left:430.2pt;margin-top:54pt;width:0;height:25.8pt;z-index:251659264;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQDhVpCU0QEAAPADAAAOAAAAZHJzL2Uyb0RvYy54bWysU02P0zAQvSPxHyzfadIgFVQ1XaEucEFQ&#10;scsP8Dp2Y2F7rLFpkn/P2GmziI89rLhMYnvezHvP493N6Cw7K4wGfMvXq5oz5SV0xp9a/u3+w6u3&#10;nMUkfCcseNXySUV+s3/5YjeErWqgB9spZFTEx+0QWt6nFLZVFWWvnIgrCMrToQZ0ItEST1WHYqDq&#10;zlZNXW+qAbALCFLFSLu38yHfl/paK5m+aB1VYrblxC2ViCU+5Fjtd2J7QhF6Iy80xDNYOGE8NV1K&#10;3Yok2A80f5RyRiJE0GklwVWgtZGqaCA16/o3NXe9CKpoIXNiWGyK/6+s/Hw+IjNdyxvOvHB0RXcJ&#10;hTn1ib1DhIEdwHuyEZA12a0hxC2BDv6Il1UMR8zSR40uf0kUG4vD0+KwGhOT86ak3dfNm82mmF89&#10;4gLG9FGBY/mn5fFCY+m/LgaL86eYqDMBr4Dc1PockzD2ve9YmgIJEZl/5ky5+bzK3Ge25S9NVs3Y&#10;r0qTB8SvKT3K9KmDRXYWNDfd9/VShTIzRBtrF1D9NOiSm2GqTOQCnBX9s9uSXTqCTwvQGQ/4t65p&#10;vFLVc/5V9aw1y36Abip3V+ygsSr+XJ5Anttf1wX++FD3PwEAAP//AwBQSwMEFAAGAAgAAAAhAOUC&#10;wwndAAAACwEAAA8AAABkcnMvZG93bnJldi54bWxMj0FPhDAQhe8m/odmTLy5raiISNmYNSboSVcP&#10;HrswC2TptKFdwH/vGA96nPe+vHmvWC92EBOOoXek4XKlQCDVrump1fDx/nSRgQjRUGMGR6jhCwOs&#10;y9OTwuSNm+kNp21sBYdQyI2GLkafSxnqDq0JK+eR2Nu70ZrI59jKZjQzh9tBJkql0pqe+ENnPG46&#10;rA/bo9UwT/ukTfzmuXp9uf08VM5XV49e6/Oz5eEeRMQl/sHwU5+rQ8mddu5ITRCDhixV14yyoTIe&#10;xcSvsmPl5i4FWRby/4byGwAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAA&#10;AAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEA&#10;AAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAOFWkJTRAQAA8AMA&#10;AA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAOUCwwndAAAA&#10;CwEAAA8AAAAAAAAAAAAAAAAAKwQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAAA1BQAA&#10;AAA=&#10;" strokecolor="black [3200]" strokeweight="2pt"><v:stroke endarrow="open"/><v:shadow on="t" color="black" opacity="24903f" origin=",.5" offset="0,.55556mm"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionRange.InsertXML($xml)
